# Commit: "fixed column charge and merged file"
#
# The file-path lookups in the mapping_column sheet (column E) pointed at
# impact_analysis/data/input/fileN.xlsx; after merging/relocating the data
# folder the impact_analysis/ prefix is no longer correct, so the sheet is
# updated to the new relative path data/input/fileN.xlsx.
#
# (Note: the canonical OOXML also shows the workbook's internal cellXfs
# style table shrinking from 7 to 6 records, with several cells' "s="
# index shifting accordingly. Every one of those records is a byte-for-byte
# duplicate of another record already in the table (no font/border/
# alignment/protection actually changes for any cell - e.g. style 4 and
# style 2 are identical <xf>s, style 3 and the new style 4 are identical,
# etc.), so there is no observable formatting change - it is pure
# bookkeeping cleanup of redundant style records performed by whatever
# tool produced that commit. That kind of raw style-table slot surgery
# has no Excel object-model surface (Range/Font/Style/Locked/... all
# resolve to the *canonical* existing record for a given set of
# properties), so it is not reproducible - or necessary - via COM and is
# intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping_column")

# --- Fix file path strings (drop the "impact_analysis/" prefix) ---
$ws.Range("E2").Value = "data/input/file1.xlsx"
$ws.Range("E5").Value = "data/input/file1.xlsx"
$ws.Range("E3").Value = "data/input/file2.xlsx"
$ws.Range("E6").Value = "data/input/file2.xlsx"
$ws.Range("E4").Value = "data/input/file3.xlsx"
$ws.Range("E7").Value = "data/input/file3.xlsx"
